$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Original")

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Stns"
$ws.Range("B1").Value = "Obs Angles"
$ws.Range("M1").Value = "Eastings"
$ws.Range("N1").Value = "Northings"

# ---------------------------------------------------------------------------
# Row 2 - Station A
# ---------------------------------------------------------------------------
$ws.Range("B2").Formula = "=(86)+(30/60)+(2/3600)"
$ws.Range("C2").Value = 187.4
$ws.Range("D2").Value = -0.0013888888888971
$ws.Range("E2").Value = 86.49916666666665
$ws.Range("F2").Value = 140.1944444444445
$ws.Range("G2").Value = 119.97
$ws.Range("H2").Value = -143.96
$ws.Range("I2").Value = -0.007200768491838459
$ws.Range("J2").Value = -0.00540057636887117
$ws.Range("K2").Value = 119.9627992315082
$ws.Range("L2").Value = -143.9654005763689
$ws.Range("M2").Value = 1000
$ws.Range("N2").Value = 1000

# ---------------------------------------------------------------------------
# Row 3 - Station B
# ---------------------------------------------------------------------------
$ws.Range("B3").Formula = "=(80)+(59/60)+(34/3600)"
$ws.Range("C3").Value = 382.7
$ws.Range("D3").Value = -0.0013888888888971
$ws.Range("E3").Value = 80.99138888888888
$ws.Range("F3").Value = 41.18583333333333
$ws.Range("G3").Value = 252.01
$ws.Range("H3").Value = 288.01
$ws.Range("I3").Value = -0.0147050912584129
$ws.Range("J3").Value = -0.011028818443794
$ws.Range("K3").Value = 251.9952949087416
$ws.Range("L3").Value = 287.9989711815562
$ws.Range("M3").Value = 1119.96
$ws.Range("N3").Value = 856.03

# ---------------------------------------------------------------------------
# Row 4 - Station C
# ---------------------------------------------------------------------------
$ws.Range("B4").Formula = "=(91)+(31/60)+(29/3600)"
$ws.Range("C4").Value = 106.1
$ws.Range("D4").Value = -0.0013888888888971
$ws.Range("E4").Value = 91.52333333333333
$ws.Range("F4").Value = 312.7091666666666
$ws.Range("G4").Value = -77.95999999999999
$ws.Range("H4").Value = 71.97
$ws.Range("I4").Value = -0.004076849183479511
$ws.Range("J4").Value = -0.003057636887605288
$ws.Range("K4").Value = -77.96407684918347
$ws.Range("L4").Value = 71.96694236311239
$ws.Range("M4").Value = 1371.96
$ws.Range("N4").Value = 1144.03

# ---------------------------------------------------------------------------
# Row 5 - Station D
# ---------------------------------------------------------------------------
$ws.Range("B5").Formula = "=(100)+(59/60)+(15/3600)"
$ws.Range("C5").Value = 364.8
$ws.Range("D5").Value = -0.0013888888888971
$ws.Range("E5").Value = 100.9861111111111
$ws.Range("F5").Value = 233.6952777777777
$ws.Range("G5").Value = -293.98
$ws.Range("H5").Value = -215.99
$ws.Range("I5").Value = -0.01401729106628959
$ws.Range("J5").Value = -0.01051296829970225
$ws.Range("K5").Value = -293.9940172910663
$ws.Range("L5").Value = -216.0005129682997
$ws.Range("M5").Value = 1294
$ws.Range("N5").Value = 1216

# ---------------------------------------------------------------------------
# Row 6 previously held Station E (A6:L6). It is replaced by three highlighted
# closing-error check cells (F6, M6, N6) that reuse the existing highlight
# fill style already present on F7 (s="2" -> fillId 2).
# ---------------------------------------------------------------------------
$ws.Range("A6:L6").ClearContents()
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("M6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("F6").Value = 140.1944444444444
$ws.Range("M6").Value = 1000.01
$ws.Range("N6").Value = 1000

# ---------------------------------------------------------------------------
# Row 7 used to hold just a lone closing-error cell (F7); it becomes the
# angle/traverse "Total" row (what used to live in row 8).
# ---------------------------------------------------------------------------
$ws.Range("F7").ClearContents()
$ws.Range("A7").Value = "Total"
$ws.Range("B7").Value = 360.0055555555556
$ws.Range("C7").Value = 1041
$ws.Range("D7").Value = -0.005555555555588398
$ws.Range("E7").Value = 360
$ws.Range("G7").Value = 0.04000000000002046
$ws.Range("H7").Value = 0.02999999999997272
$ws.Range("I7").Value = -0.04000000000002046
$ws.Range("J7").Value = -0.02999999999997271
$ws.Range("K7").Value = -0
$ws.Range("L7").Value = 0

# ---------------------------------------------------------------------------
# Row 8 used to hold the "Total" row; it becomes the fractional-misclosure
# caption row.
# ---------------------------------------------------------------------------
$ws.Range("C8:L8").ClearContents()
$ws.Range("A8").Value = "Fractional Misclosure"
$ws.Range("B8").Value = "1 / 20820"

# ---------------------------------------------------------------------------
# Rows 9 and 10 are no longer present at all.
# ---------------------------------------------------------------------------
$ws.Range("A9:N10").Delete() | Out-Null

# ---------------------------------------------------------------------------
# Column widths: only column B keeps an explicit width now.
# ---------------------------------------------------------------------------
$ws.Range("A:A").ColumnWidth = 8
$ws.Range("C:C").ColumnWidth = 8
$ws.Range("D:D").ColumnWidth = 8
$ws.Range("E:E").ColumnWidth = 8
$ws.Range("F:F").ColumnWidth = 8
$ws.Range("B:B").ColumnWidth = 14.28515625

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("C6").Select() | Out-Null
